# Final push for project submission.
# Marks "In-game controls screen" and "Menus" as Implemented = Yes
# (previously "No" / "WIP"), and adds a confident remark in the Location
# column for the "Fun, funny, interesting, surprising" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 - "In-game controls screen": No -> Yes
$ws.Range("B34").Value = "Yes"
$ws.Range("B34").Interior.Color = 5296274

# Row 35 - "Menus": WIP -> Yes
$ws.Range("B35").Value = "Yes"
$ws.Range("B35").Interior.Color = 5296274

# Row 37 - "Fun, funny, interesting, surprising": add a location/notes remark
$ws.Range("C37").Value = "Boy do I hope so!"

# Move the active selection to reflect where the author was last working
$ws.Range("C34").Select()
